$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching style of existing header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for new columns I (I0) and J (IF)
$iValues = @(8,7,9,7,5,9,6,6,7,2,7,6,7,6,6)
$jValues = @(8,7,9,9,5,9,6,6,7,2,7,6,7,6,6)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
